$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Device")

# Row 4 previously held "Device 3" data; it is replaced with the data that used
# to be "Device 4" (row 5), i.e. the "Device 3" row was removed and everything
# shifted up by one row.
$ws.Range("A4").Value = "Device 4"
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 10
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = "0x2a"
$ws.Range("G4").Value = "Temperature "
$ws.Range("H4").Value = "Celsius"

# Row 5 previously held "Device 4"; now holds what used to be "Device 5" (row 6)
$ws.Range("A5").Value = "Device 5"
$ws.Range("E5").Value = 5

# Row 6 previously held "Device 5"; now holds what used to be "Device 6" (row 7)
$ws.Range("A6").Value = "Device 6"
$ws.Range("E6").Value = 6

# Row 7 previously held "Device 6"; now becomes a new "Device 7" entry
$ws.Range("A7").Value = "Device 7"
$ws.Range("E7").Value = 7

# Update the selection to match the new active cell/range (the whole row 4
# was selected, e.g. after removing the old "Device 3" row).
$ws.Range("A4:XFD4").Select()
